# feat: add 2022-Q3 data
#
# The workbook gains a new quarterly sheet "2022-Q3", inserted right after
# the "总计" (totals) sheet and before "2022-Q1" (so every existing quarter
# sheet shifts one position to the right). The "总计" summary sheet gets a
# new row for every quarter (including a brand-new trailing row that
# reproduces the former last row, since the sheet count grew by one) and its
# existing rows shift down to make room for the new 2022-Q3 entry at the
# top of the data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet: shift the quarter rows down and append a new one.
# ---------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

# Copy the formatting of the last existing data row (A5) down into the new
# A6 so the new row's style (bold/bordered like the other "A" column cells)
# matches the rest of the column before we fill in its values.
$zj.Range("A5").Copy()
$zj.Range("A6").PasteSpecial(-4122)

$zj.Range("B2").Value = "2022-Q3"
$zj.Range("C2").Value = 3
$zj.Range("D2").Value = 0.19

$zj.Range("B3").Value = "2022-Q1"
$zj.Range("C3").Value = 3
$zj.Range("D3").Value = 0.2

$zj.Range("B4").Value = "2021-Q4"
$zj.Range("C4").Value = 3
$zj.Range("D4").Value = 0.25

$zj.Range("B5").Value = "2021-Q3"
$zj.Range("C5").Value = 3
$zj.Range("D5").Value = 0.23

$zj.Range("A6").Value = 4
$zj.Range("B6").Value = "2021-Q2"
$zj.Range("C6").Value = 3
$zj.Range("D6").Value = 0.23

# ---------------------------------------------------------------------
# 2) Create the new "2022-Q3" sheet by duplicating "2022-Q1" (same column
#    layout/styling) and placing the copy immediately before it, then
#    overwrite its data with the 2022-Q3 figures.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Copy($q1, $null)
$q3 = $wb.Worksheets.Item("2022-Q1 (2)")
$q3.Name = "2022-Q3"

# The B (fund code, e.g. "000593" has a meaningful leading zero) and D:G
# columns hold text-looking numbers (e.g. "3.41", "0.0631") that must stay
# text (not be coerced to numeric, which would drop the leading zero or
# trailing decimal zeros), so mark that range as Text before writing values.
$q3.Range("B2:B4").NumberFormat = "@"
$q3.Range("D2:G4").NumberFormat = "@"

$q3.Range("B2").Value = "118002"
$q3.Range("C2").Value = "易方达标普全球高端消费品指数增强A（QDII）人民币"
$q3.Range("D2").Value = "1.85"
$q3.Range("E2").Value = "93.04"
$q3.Range("F2").Value = "3.41"
$q3.Range("G2").Value = "0.0631"
$q3.Range("H2").Value = 10

$q3.Range("B3").Value = "000593"
$q3.Range("C3").Value = "易方达标普全球高端消费品指数增强（QDII）美元现汇"
$q3.Range("D3").Value = "1.85"
$q3.Range("E3").Value = "93.04"
$q3.Range("F3").Value = "3.41"
$q3.Range("G3").Value = "0.0631"
$q3.Range("H3").Value = 10

$q3.Range("B4").Value = "005676"
$q3.Range("C4").Value = "易方达标普全球高端消费品指数增强C（QDII）人民币"
$q3.Range("D4").Value = "1.85"
$q3.Range("E4").Value = "93.04"
$q3.Range("F4").Value = "3.41"
$q3.Range("G4").Value = "0.0631"
$q3.Range("H4").Value = 10

# The rest of the existing quarter sheets ("2022-Q1", "2021-Q4", "2021-Q3",
# "2021-Q2") keep their original data and formatting untouched; only their
# tab position shifts right, which naturally happened by inserting the new
# "2022-Q3" sheet before "2022-Q1" above.

$zj.Select()
